$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.817.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.932.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '376.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0852'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.408.13'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("B15").Value = 'Uniswap'
$ws.Range("C15").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '12.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +68.41%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.939.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.991'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '50.807.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("E20").Value = '  -6.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.108'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("E36").Value = '  -2.57%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +4.01%  '
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '123.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("E45").Value = '  +6.36%  '
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.001.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  -5.50%  '
$ws.Range("E50").Value = '  -5.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.17%  '
